$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-7 from 45175 (2023-09-06)
# to 45183 (2023-09-14), matching the automatic file update.
$ws.Range("C2:C7").Value = 45183
